$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nk")

$words = @("Porca", "Gay", "Arcobaleno", "Dio", "Comodinodio", "Start", "Pisellino", "Porco spino")

$row = 11
foreach ($word in $words) {
    $ws.Cells.Item($row, 1).Value = $word
    $row = $row + 1
}
